$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 64
$ws.Range("H64").Value = 7491
$ws.Range("J64").Value = 8500
$ws.Range("L64").Value = 8500
$ws.Range("N64").Value = -8996
# Row 67
$ws.Range("H67").Value = 7491
$ws.Range("J67").Value = 8500
$ws.Range("L67").Value = 8500
$ws.Range("N67").Value = -10216
# Row 70
$ws.Range("H70").Value = 9617627
$ws.Range("J70").Value = 11366059
$ws.Range("L70").Value = 34098177
$ws.Range("N70").Value = -34098717
# Row 73
$ws.Range("H73").Value = 9617627
$ws.Range("J73").Value = 11366059
$ws.Range("L73").Value = 34098177
$ws.Range("N73").Value = -34100049
# Row 96
$ws.Range("H96").Value = 2199.8572
$ws.Range("I96").Value = 513.25
$ws.Range("K96").Value = 1539.75
$ws.Range("M96").Value = -166.75
# Row 100
$ws.Range("H100").Value = 3308.9285
$ws.Range("I100").Value = 2124.3333
$ws.Range("K100").Value = 2124.3333
$ws.Range("M100").Value = -1583.3333

$ws = $wb.Worksheets.Item("ARM")
# Row 122
$ws.Range("H122").Value = 4986.5293
$ws.Range("I122").Value = 3608.0476
$ws.Range("K122").Value = 10824.1428
$ws.Range("M122").Value = -8374.1428
# Row 126
$ws.Range("H126").Value = 9997.5
$ws.Range("I126").Value = 9997.5
$ws.Range("K126").Value = 29992.5
$ws.Range("M126").Value = -27522.5

$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 3035.6365
$ws.Range("I20").Value = 2355.4285
$ws.Range("K20").Value = 2355.4285
$ws.Range("M20").Value = -2108.4285
# Row 86
$ws.Range("H86").Value = 1599.3334
$ws.Range("I86").Value = 798
$ws.Range("J86").Value = 2000
$ws.Range("K86").Value = 798
$ws.Range("L86").Value = 2000
$ws.Range("M86").Value = 325
$ws.Range("N86").Value = -4246
# Row 89
$ws.Range("H89").Value = 1599.3334
$ws.Range("I89").Value = 798
$ws.Range("J89").Value = 2000
$ws.Range("K89").Value = 3990
$ws.Range("L89").Value = 10000
$ws.Range("M89").Value = 1626
$ws.Range("N89").Value = -21232
# Row 105
$ws.Range("H105").Value = 2381.7334
$ws.Range("I105").Value = 2048.24
$ws.Range("K105").Value = 2048.24
$ws.Range("M105").Value = -301.2399999999998
# Row 107
$ws.Range("H107").Value = 1529.775
$ws.Range("I107").Value = 1399.7241
$ws.Range("J107").Value = 1872.6364
$ws.Range("K107").Value = 1399.7241
$ws.Range("L107").Value = 1872.6364
$ws.Range("M107").Value = 520.2759000000001
$ws.Range("N107").Value = -5712.6364
# Row 134
$ws.Range("H134").Value = 2754.8572
$ws.Range("I134").Value = 2351
$ws.Range("J134").Value = 2916.4
$ws.Range("K134").Value = 7053
$ws.Range("L134").Value = 8749.200000000001
$ws.Range("M134").Value = -4518
$ws.Range("N134").Value = -13819.2

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 3864.3489
$ws.Range("I31").Value = 2344.0715
$ws.Range("J31").Value = 4598.276
$ws.Range("K31").Value = 2344.0715
$ws.Range("L31").Value = 4598.276
$ws.Range("M31").Value = -2049.0715
$ws.Range("N31").Value = -5188.276
# Row 34
$ws.Range("H34").Value = 3864.3489
$ws.Range("I34").Value = 2344.0715
$ws.Range("J34").Value = 4598.276
$ws.Range("K34").Value = 2344.0715
$ws.Range("L34").Value = 4598.276
$ws.Range("M34").Value = -2142.0715
$ws.Range("N34").Value = -5002.276
# Row 56
$ws.Range("H56").Value = 0
$ws.Range("J56").Value = 0
$ws.Range("L56").Value = 0
$ws.Range("N56").Value = ""

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 459.81818
$ws.Range("I5").Value = 417.66666
$ws.Range("J5").Value = 649.5
$ws.Range("K5").Value = 1252.99998
$ws.Range("L5").Value = 1948.5
$ws.Range("M5").Value = -1140.99998
$ws.Range("N5").Value = -2172.5
# Row 26
$ws.Range("H26").Value = 1733.3334
$ws.Range("I26").Value = 2000
$ws.Range("K26").Value = 6000
$ws.Range("M26").Value = -5712
# Row 61
$ws.Range("H61").Value = 576.8333
$ws.Range("I61").Value = 152
$ws.Range("J61").Value = 789.25
$ws.Range("K61").Value = 456
$ws.Range("L61").Value = 2367.75
$ws.Range("M61").Value = -241
$ws.Range("N61").Value = -2797.75
# Row 113
$ws.Range("H113").Value = 1482
$ws.Range("J113").Value = 1405.9333
$ws.Range("L113").Value = 4217.7999
$ws.Range("N113").Value = -8557.7999
# Row 135
$ws.Range("H135").Value = 459.81818
$ws.Range("I135").Value = 417.66666
$ws.Range("J135").Value = 649.5
$ws.Range("K135").Value = 3758.99994
$ws.Range("L135").Value = 5845.5
$ws.Range("M135").Value = -1223.99994
$ws.Range("N135").Value = -10915.5
# Row 136
$ws.Range("H136").Value = 1407.375
$ws.Range("I136").Value = 1407.375
$ws.Range("K136").Value = 4222.125
$ws.Range("M136").Value = 877.875

$ws = $wb.Worksheets.Item("GSM")
# Row 63
$ws.Range("H63").Value = 64683
$ws.Range("J63").Value = 64683
$ws.Range("L63").Value = 64683
$ws.Range("N63").Value = -66055
# Row 66
$ws.Range("H66").Value = 64683
$ws.Range("J66").Value = 64683
$ws.Range("L66").Value = 194049
$ws.Range("N66").Value = -200913
# Row 70
$ws.Range("H70").Value = 4490.2197
$ws.Range("I70").Value = 4433
$ws.Range("K70").Value = 4433
$ws.Range("M70").Value = -4163
# Row 73
$ws.Range("H73").Value = 4490.2197
$ws.Range("I73").Value = 4433
$ws.Range("K73").Value = 4433
$ws.Range("M73").Value = -3497
# Row 80
$ws.Range("H80").Value = 3139.8
$ws.Range("I80").Value = 2844.1667
$ws.Range("K80").Value = 2844.1667
$ws.Range("M80").Value = -1846.1667
# Row 83
$ws.Range("H83").Value = 3139.8
$ws.Range("I83").Value = 2844.1667
$ws.Range("K83").Value = 14220.8335
$ws.Range("M83").Value = -9228.833500000001
# Row 132
$ws.Range("H132").Value = 3331.0625
$ws.Range("I132").Value = 3116.5
$ws.Range("K132").Value = 9349.5
$ws.Range("M132").Value = -6819.5

$ws = $wb.Worksheets.Item("LTW")
# Row 5
$ws.Range("H5").Value = 23254.75
$ws.Range("I5").Value = 12004.5
$ws.Range("J5").Value = 34505
$ws.Range("K5").Value = 12004.5
$ws.Range("L5").Value = 34505
$ws.Range("M5").Value = -11891.5
$ws.Range("N5").Value = -34731
# Row 21
$ws.Range("H21").Value = 20006
$ws.Range("I21").Value = 20006
$ws.Range("K21").Value = 20006
$ws.Range("M21").Value = -19832
# Row 40
$ws.Range("H40").Value = 5486.2856
$ws.Range("I40").Value = 6412.7144
$ws.Range("J40").Value = 4559.857
$ws.Range("K40").Value = 6412.7144
$ws.Range("L40").Value = 4559.857
$ws.Range("M40").Value = -6276.7144
$ws.Range("N40").Value = -4831.857
# Row 82
$ws.Range("H82").Value = 2718.7778
$ws.Range("I82").Value = 7000
$ws.Range("J82").Value = 2183.625
$ws.Range("K82").Value = 7000
$ws.Range("L82").Value = 2183.625
$ws.Range("M82").Value = -6639
$ws.Range("N82").Value = -2905.625
# Row 85
$ws.Range("H85").Value = 2718.7778
$ws.Range("I85").Value = 7000
$ws.Range("J85").Value = 2183.625
$ws.Range("K85").Value = 7000
$ws.Range("L85").Value = 2183.625
$ws.Range("M85").Value = -5752
$ws.Range("N85").Value = -4679.625
# Row 122
$ws.Range("H122").Value = 23026.834
$ws.Range("I122").Value = 23950.8
$ws.Range("K122").Value = 71852.39999999999
$ws.Range("M122").Value = -69402.39999999999
# Row 132
$ws.Range("H132").Value = 4000
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").Value = ""

$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Range("H122").Value = 9965.429
$ws.Range("I122").Value = 3668
$ws.Range("K122").Value = 11004
$ws.Range("M122").Value = -8554
# Row 132
$ws.Range("H132").Value = 3696.7727
$ws.Range("I132").Value = 3464.4211
$ws.Range("K132").Value = 10393.2633
$ws.Range("M132").Value = -7863.263300000001
# Row 136
$ws.Range("H136").Value = 3873.0625
$ws.Range("I136").Value = 3366.6667
$ws.Range("K136").Value = 10100.0001
$ws.Range("M136").Value = -7550.000100000001
